# daily auto push: 2026-01-10 18:43 UTC
# Insert two new rows of data at the top of the rolling log (row 623),
# pushing the existing rows down by 2. The worksheet dimension and the
# inserted rows pick up the two newest readings for 2026/01/10 and
# 2026/01/11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (formerly rows 623:664) down by two rows.
$ws.Rows("623:624").Insert()

# Pre-format column A for the new rows as Text so the "yyyy/mm/dd"
# looking strings are stored as literal text (matching the rest of the
# column) instead of being auto-parsed into date serial numbers.
$ws.Range("A623:A624").NumberFormat = "@"

# New row 623: 2026/01/10 (Sat)
$ws.Range("A623").Value = "2026/01/10"
$ws.Range("B623").Value = "土"
$ws.Range("C623").Value = 23
$ws.Range("D623").Value = 200

# New row 624: 2026/01/11 (Sun)
$ws.Range("A624").Value = "2026/01/11"
$ws.Range("B624").Value = "日"
$ws.Range("C624").Value = 2
$ws.Range("D624").Value = 146

# Strip the formatting we applied above so the new cells end up with no
# explicit style, consistent with the rest of the sheet.
$ws.Range("A623:D624").ClearFormats()
